$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing input values (B4/B5/B6/B9) and dependent formula results recalc automatically ---
$ws.Range("B4").Value = 40
$ws.Range("B5").Value = 50
$ws.Range("B6").Value = 60
$ws.Range("B9").Value = 120

# --- New header label ---
$ws.Range("G2").Value = "Verification tests:"

# --- Verification test block 1: Inputs (rows 4-7) ---
$ws.Range("F4").Value = "Inputs"
$ws.Range("G4").Value = "carb"
$ws.Range("H4").Value = 40
$ws.Range("J4").Value = "carb"
$ws.Range("K4").Value = 40
$ws.Range("M4").Value = "carb"
$ws.Range("N4").Value = 0

$ws.Range("G5").Value = "fat"
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = "fat"
$ws.Range("K5").Value = 50
$ws.Range("M5").Value = "fat"
$ws.Range("N5").Value = 50

$ws.Range("G6").Value = "protein"
$ws.Range("H6").Value = 60
$ws.Range("J6").Value = "protein"
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = "protein"
$ws.Range("N6").Value = 60

$ws.Range("G7").Value = "ratio"
$ws.Range("H7").Value = 80
$ws.Range("J7").Value = "ratio"
$ws.Range("K7").Value = 100
$ws.Range("M7").Value = "ratio"
$ws.Range("N7").Value = 120

# --- Verification test block 1: Expected (rows 9-12) ---
$ws.Range("F9").Value = "Expected "
$ws.Range("G9").Value = "bolus"
$ws.Range("H9").Value = 40
$ws.Range("J9").Value = "bolus"
$ws.Range("K9").Value = 40
$ws.Range("M9").Value = "bolus"
$ws.Range("N9").Value = 0

$ws.Range("G10").Value = "square"
$ws.Range("H10").Value = 30
$ws.Range("J10").Value = "square"
$ws.Range("K10").Value = 45
$ws.Range("M10").Value = "square"
$ws.Range("N10").Value = 29

$ws.Range("G11").Value = "duration"
$ws.Range("H11").Value = 8
$ws.Range("J11").Value = "duration"
$ws.Range("K11").Value = 9.5
$ws.Range("M11").Value = "duration"
$ws.Range("N11").Value = 10.8

$ws.Range("G12").Value = "passed"
$ws.Range("J12").Value = "passed"
$ws.Range("M12").Value = "passed"

# --- New notes under row 14 ---
$ws.Range("A16").Value = "Note - duration should be entered as 4 hrs minimum for manual pump"
$ws.Range("A17").Value = "or 7 hours minimum for Loop or OpenAPS."

# --- Verification test block 2: Inputs (rows 15-18) ---
$ws.Range("F15").Value = "Inputs"
$ws.Range("G15").Value = "carb"
$ws.Range("H15").Value = 40
$ws.Range("J15").Value = "carb"
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = "carb"
$ws.Range("N15").Value = 40

$ws.Range("G16").Value = "fat"
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = "fat"
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = "fat"
$ws.Range("N16").Value = 50

$ws.Range("G17").Value = "protein"
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = "protein"
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = "protein"
$ws.Range("N17").Value = 60

$ws.Range("G18").Value = "ratio"
$ws.Range("H18").Value = 80
$ws.Range("J18").Value = "ratio"
$ws.Range("K18").Value = 100
$ws.Range("M18").Value = "ratio"
$ws.Range("N18").Value = 120

# --- Verification test block 2: Expected (rows 20-23) ---
$ws.Range("F20").Value = "Expected "
$ws.Range("G20").Value = "bolus"
$ws.Range("H20").Value = 40
$ws.Range("J20").Value = "bolus"
$ws.Range("K20").Value = "NA"
$ws.Range("M20").Value = "bolus"
$ws.Range("N20").Value = 40

$ws.Range("G21").Value = "square"
$ws.Range("H21").Value = "NA"
$ws.Range("J21").Value = "square"
$ws.Range("K21").Value = "NA"
$ws.Range("M21").Value = "square"
$ws.Range("N21").Value = 58

$ws.Range("G22").Value = "duration"
$ws.Range("H22").Value = "NA"
$ws.Range("J22").Value = "duration"
$ws.Range("K22").Value = "NA"
$ws.Range("M22").Value = "duration"
$ws.Range("N22").Value = 10.8

$ws.Range("G23").Value = "passed"
$ws.Range("J23").Value = "passed"
$ws.Range("M23").Value = "passed"

# --- Styling ---
# "NA" cells -> right aligned
$naCells = @("K20","H21","K21","H22","K22")
foreach ($ref in $naCells) {
    $ws.Range($ref).HorizontalAlignment = -4152
}

# "passed" result cells -> green font
$passedCells = @("G12","J12","M12","G23","J23","M23")
foreach ($ref in $passedCells) {
    $ws.Range($ref).Font.Color = 32768
}

# Header labels ("Verification tests:", "Inputs", "Expected ") -> blue font
$headerCells = @("G2","F4","F9","F15","F20")
foreach ($ref in $headerCells) {
    $ws.Range($ref).Font.Color = 16711680
}

# --- Selection matches the target view state ---
$ws.Range("G27").Select()
